# Deploy the implementation guide: refresh the generated Metadata sheet
# (Status + Date) to reflect the new publication run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $label = $ws.Cells.Item($r, 1).Text

    if ($label -eq "Status") {
        $ws.Cells.Item($r, 2).Value = "draft"
    }
    elseif ($label -eq "Date") {
        $ws.Cells.Item($r, 2).Value = "2023-08-01T16:12:28+00:00"
    }
}
